$wb = $excel.ActiveWorkbook

# Update the Users sheet: the user name changed from "Nicole Bicho" to "Drew Koecher"
$usersSheet = $wb.Worksheets.Item("Users")
$usersSheet.Range("A2").Value = "Drew Koecher"

# Make the Users sheet the active tab, with cell C8 selected
$usersSheet.Activate()
$usersSheet.Range("C8").Select()
